$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has columns A..H:
#   A=CUSTOMER.ID  B=DEAL.DATE  C=VALUE.DATE  D=MATURITY.DATE  E=CURRENCY
#   F=PRINCIPAL    G=INTEREST.RATE  H=AUTO.ROLLOVER
#
# Insert 7 new columns before H so AUTO.ROLLOVER (and its value) shifts to
# column O, then add the new header labels in H..N plus two more headers
# after AUTO.ROLLOVER (now P, Q). Content/formatting of the existing
# columns (including the AUTO.ROLLOVER column itself) is carried along by
# the column insert.
$ws.Range("H1:N1").EntireColumn.Insert()

$ws.Range("H1").Value = "INTEND.DATE"
$ws.Range("I1").Value = "CUST.REMARKS:1"
$ws.Range("J1").Value = "TAX.INTEREST.TYPE:1"
$ws.Range("K1").Value = "DRAWDOWN.ACCOUNT"
$ws.Range("L1").Value = "PRIN.LIQ.ACCT"
$ws.Range("M1").Value = "INT.LIQ.ACCT"
$ws.Range("N1").Value = "CHRG.LIQ.ACCT"

$ws.Range("P1").Value = "FINAL.MATURITY"
$ws.Range("Q1").Value = "EXP.DATE"

# Match the column width used by the rest of the data columns (e.g. G) for
# the newly inserted columns.
$ws.Range("H1:N1").EntireColumn.ColumnWidth = 13.5

# Scroll the sheet a bit and leave the selection on the new CUST.REMARKS:1
# column like the author did while reviewing the new fields.
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("I13").Select()
